$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (matches the diff: rows 245-247)
$data = @(
    @(245, 44319, 1, 2, 43.62050163576881),
    @(246, 44320, 0, 2, 43.62050163576881),
    @(247, 44321, 0, 2, 43.62050163576881)
)

foreach ($row in $data) {
    $r = $row[0]

    # Copy formatting (not values) from the row above so the new row's
    # style indices exactly match the existing data (date cell style, etc.)
    $ws.Range("A" + ($r - 1) + ":D" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
